# Atualização automática de ITAQUI.xlsx
#  - Renomeia "Paineis DARQ" -> "PAINEIS DARQ"
#  - Renomeia "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
#  - Remove a aba "Desarquivamentos Pendentes"

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete() | Out-Null
$excel.DisplayAlerts = $true
